$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")
$ws.Activate()

# Row 6 - Task 1 (Jogar o jogo...)
$ws.Range("D6").Value = 8
$ws.Range("K6").Value = 2

# Row 7 - Task 2
$ws.Range("K7").Value = 1

# Row 11 - Task 6: text change "Meeting semanal" -> "Meetings"
$ws.Range("C11").Value = "Meetings"
$ws.Range("D11").Value = 3.5
$ws.Range("K11").Value = 1.5

# Row 13 - Task 8
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 1

# Row 14 - Task 9
$ws.Range("K14").Value = 1

# Row 15 - Task 10: new text + values
$ws.Range("C15").Value = "Fazer pdf dos 3 User Stories mais votados pela equipa e submeter no moodle"
$ws.Range("D15").Value = 1
$ws.Range("K15").Value = 1

# Recalculate so dependent formulas (sums, ideal burndown) and the chart
# series that read them are refreshed
$excel.Calculate()
$co = $ws.ChartObjects().Item(1)
$co.Chart.Refresh()

# Sheet view: scroll position + selection (mirrors the author scrolling
# down to row 7 and leaving the cursor on U4)
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("U4").Select()
